$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "307.45"
    "E2" = "-1.17%"
    "D3" = "37.17"
    "E3" = "-1.21%"
    "D4" = "5.114"
    "E4" = "0.28%"
    "D5" = "0.07785"
    "E5" = "0.10%"
    "D6" = "8.249"
    "E6" = "0.44%"
    "D7" = "1.882"
    "E7" = "-0.63%"
    "E8" = "4.16%"
    "D9" = "0.9261"
    "E9" = "0.43%"
    "D10" = "0.1094"
    "E10" = "-10.15%"
    "D11" = "0.1897"
    "E11" = "-0.69%"
    "D12" = "0.08887"
    "E12" = "-3.33%"
    "D13" = "0.03354"
    "E13" = "-2.18%"
    "D14" = "0.09577"
    "E14" = "-1.09%"
    "D15" = "0.001380"
    "E15" = "0.59%"
    "D16" = "0.005732"
    "E16" = "-3.00%"
    "D17" = "3.516"
    "E17" = "-1.19%"
    "D18" = "4.420"
    "E18" = "1.00%"
    "D19" = "0.3367"
    "E19" = "-1.05%"
    "D20" = "6.280"
    "E20" = "19.63%"
    "D21" = "0.1277"
    "E21" = "-1.51%"
    "D22" = "0.2503"
    "E22" = "-3.42%"
    "D23" = "0.04383"
    "E23" = "0.54%"
    "D24" = "0.001191"
    "E24" = "-1.80%"
    "D25" = "0.004253"
    "E25" = "-0.07%"
    "D26" = "0.0001303"
    "E26" = "0.16%"
    "D39" = "0.02159"
    "E39" = "3.17%"
    "D40" = "0.04995"
    "E40" = "-3.73%"
    "D41" = "0.007463"
    "E41" = "-3.12%"
    "D42" = "0.1350"
    "E42" = "0.35%"
    "D43" = "0.008650"
    "E43" = "-11.18%"
    "D44" = "0.002039"
    "E44" = "-1.08%"
    "D45" = "0.008024"
    "E45" = "-9.93%"
    "D46" = "0.00006166"
    "E46" = "-7.72%"
    "E47" = "0.16%"
    "D48" = "0.003288"
    "E48" = "11.95%"
    "D49" = "0.001002"
    "E49" = "-16.54%"
    "D50" = "0.00002105"
    "E50" = "0.16%"
    "D51" = "0.0002005"
    "E51" = "0.16%"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = "'" + $updates[$cell]
    $ws.Range($cell).Style = "Normal"
}
